$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.544.38'
$ws.Range("E2").Value = '  +3.24%  '
$ws.Range("D3").Value = '3.347.01'
$ws.Range("E3").Value = '  +4.10%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '193.18'
$ws.Range("E5").Value = '  +5.95%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '591.98'
$ws.Range("E6").Value = '  +2.44%  '
$ws.Range("E8").Value = '  +1.39%  '
$ws.Range("E9").Value = '  +4.07%  '
$ws.Range("E10").Value = '  +2.79%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.422'
$ws.Range("E11").Value = '  +2.66%  '
$ws.Range("D12").Value = '3.926.67'
$ws.Range("E12").Value = '  +3.98%  '
$ws.Range("E13").Value = '  +1.24%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.35'
$ws.Range("E14").Value = '  +3.19%  '
$ws.Range("D15").Value = '69.522.86'
$ws.Range("E15").Value = '  +3.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000172'
$ws.Range("E16").Value = '  +2.43%  '
$ws.Range("D17").Value = '3.348.19'
$ws.Range("E17").Value = '  +4.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.82'
$ws.Range("E18").Value = '  +1.83%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '438.79'
$ws.Range("E19").Value = '  +12.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.71'
$ws.Range("E20").Value = '  +2.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.77'
$ws.Range("E21").Value = '  +3.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.46'
$ws.Range("E22").Value = '  +4.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("D24").Value = '3.489.24'
$ws.Range("E24").Value = '  +3.90%  '
$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000122'
$ws.Range("E25").Value = '  +4.61%  '
$ws.Range("B26").Value = 'Polygon'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.517'
$ws.Range("E26").Value = '  +1.37%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.193'
$ws.Range("E27").Value = '  +4.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.57'
$ws.Range("E28").Value = '  +0.79%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("E30").Value = '  +2.83%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '23.09'
$ws.Range("E31").Value = '  +2.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.59'
$ws.Range("E32").Value = '  +1.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.04'
$ws.Range("E33").Value = '  +1.66%  '
$ws.Range("E34").Value = '  +2.70%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("E36").Value = '  +3.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '164.37'
$ws.Range("E37").Value = '  +1.89%  '
$ws.Range("E38").Value = '  +2.25%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '27.13'
$ws.Range("E39").Value = '  +3.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.811'
$ws.Range("E40").Value = '  +1.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.57'
$ws.Range("E41").Value = '  +1.36%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '2.757.17'
$ws.Range("E42").Value = '  +6.23%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.50'
$ws.Range("E43").Value = '  +1.13%  '
$ws.Range("E44").Value = '  +3.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0689'
$ws.Range("E45").Value = '  +1.57%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.33'
$ws.Range("E46").Value = '  +2.81%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '343.76'
$ws.Range("E47").Value = '  +3.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.68'
$ws.Range("E48").Value = '  +0.02%  '
$ws.Range("E49").Value = '  +3.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '32.65'
$ws.Range("E50").Value = '  +6.61%  '
$ws.Range("E51").Value = '  +5.19%  '
